# Weekly driver report update for 2025-04-21
# Updates the "Bad Drivers" and "Good Drivers" tables on the Driver Summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextLiteral($cell, $text) {
    # Force a plain text value (avoids Excel auto-converting date-shaped
    # strings like "2024-07-23" into a date serial number).
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# ---- Bad Drivers table (rows 3-7) ----
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 692
$ws.Cells.Item(3, 4).Value = 95

$ws.Cells.Item(4, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.80.0.7"
$ws.Cells.Item(4, 2).Value = 31
$ws.Cells.Item(4, 3).Value = 1377
$ws.Cells.Item(4, 4).Value = 98.7

$ws.Cells.Item(5, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.18.2"
$ws.Cells.Item(5, 2).Value = 21
$ws.Cells.Item(5, 3).Value = 1857
$ws.Cells.Item(5, 4).Value = 98.8

$ws.Cells.Item(6, 3).Value = 1108

$ws.Cells.Item(7, 2).Value = 77
$ws.Cells.Item(7, 3).Value = 5034

# ---- Good Drivers table (rows 15-27 updated in place) ----
$ws.Cells.Item(15, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3"
$ws.Cells.Item(15, 2).Value = 18721
$ws.Cells.Item(15, 4).Value = 99.90000000000001
Set-TextLiteral $ws.Cells.Item(15, 5) "2024-07-23"

$ws.Cells.Item(16, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1"
$ws.Cells.Item(16, 2).Value = 69578
$ws.Cells.Item(16, 4).Value = 99.90000000000001
Set-TextLiteral $ws.Cells.Item(16, 5) "2023-08-14"

$ws.Cells.Item(17, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8"
$ws.Cells.Item(17, 2).Value = 331283
Set-TextLiteral $ws.Cells.Item(17, 5) "2023-05-08"

$ws.Cells.Item(18, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.200.0.6"
$ws.Cells.Item(18, 2).Value = 143808
Set-TextLiteral $ws.Cells.Item(18, 5) "2023-01-16"

$ws.Cells.Item(19, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.190.0.4"
$ws.Cells.Item(19, 2).Value = 287148
$ws.Cells.Item(19, 4).Value = 99.90000000000001
Set-TextLiteral $ws.Cells.Item(19, 5) "2022-11-22"

$ws.Cells.Item(20, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4"
$ws.Cells.Item(20, 2).Value = 96526
$ws.Cells.Item(20, 4).Value = 99.90000000000001
Set-TextLiteral $ws.Cells.Item(20, 5) "2022-08-13"

$ws.Cells.Item(21, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 22.30.0.11"
$ws.Cells.Item(21, 2).Value = 170510
Set-TextLiteral $ws.Cells.Item(21, 5) "2021-01-19"

$ws.Cells.Item(22, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.30.0.11"
$ws.Cells.Item(22, 2).Value = 67111
Set-TextLiteral $ws.Cells.Item(22, 5) "2021-01-19"

$ws.Cells.Item(23, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.10.0.7"
$ws.Cells.Item(23, 2).Value = 66577
Set-TextLiteral $ws.Cells.Item(23, 5) "2020-10-19"

$ws.Cells.Item(24, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.0.1.1"
$ws.Cells.Item(24, 2).Value = 15734
$ws.Cells.Item(24, 4).Value = 99.90000000000001
Set-TextLiteral $ws.Cells.Item(24, 5) "2020-09-28"

$ws.Cells.Item(25, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 22.0.1.1"
$ws.Cells.Item(25, 2).Value = 52096
$ws.Cells.Item(25, 4).Value = 100
Set-TextLiteral $ws.Cells.Item(25, 5) "2020-09-28"

$ws.Cells.Item(26, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.60.2.1"
$ws.Cells.Item(26, 2).Value = 26241
$ws.Cells.Item(26, 4).Value = 100
Set-TextLiteral $ws.Cells.Item(26, 5) "2019-12-14"

$ws.Cells.Item(27, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.11.3"
$ws.Cells.Item(27, 2).Value = 161874
$ws.Cells.Item(27, 4).Value = 100
Set-TextLiteral $ws.Cells.Item(27, 5) "2019-09-05"

# ---- Good Drivers table: new rows 28-37 (copy formatting from row 27) ----
$ws.Range("A27:E27").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)
$ws.Cells.Item(28, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.40.2.2"
$ws.Cells.Item(28, 2).Value = 88435
$ws.Cells.Item(28, 4).Value = 99.90000000000001
Set-TextLiteral $ws.Cells.Item(28, 5) "2019-08-31"

$ws.Range("A27:E27").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)
$ws.Cells.Item(29, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5"
$ws.Cells.Item(29, 2).Value = 144782
$ws.Cells.Item(29, 4).Value = 99.90000000000001
Set-TextLiteral $ws.Cells.Item(29, 5) "2019-08-25"

$ws.Range("A27:E27").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)
$ws.Cells.Item(30, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.30.4.1"
$ws.Cells.Item(30, 2).Value = 13016
$ws.Cells.Item(30, 4).Value = 100
Set-TextLiteral $ws.Cells.Item(30, 5) "2019-07-29"

$ws.Range("A27:E27").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)
$ws.Cells.Item(31, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.10.2"
$ws.Cells.Item(31, 2).Value = 20227
$ws.Cells.Item(31, 4).Value = 100
Set-TextLiteral $ws.Cells.Item(31, 5) "2019-05-11"

$ws.Range("A27:E27").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122)
$ws.Cells.Item(32, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.9.1"
$ws.Cells.Item(32, 2).Value = 34065
$ws.Cells.Item(32, 4).Value = 100
Set-TextLiteral $ws.Cells.Item(32, 5) "2019-04-28"

$ws.Range("A27:E27").Copy()
$ws.Range("A33:E33").PasteSpecial(-4122)
$ws.Cells.Item(33, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.10.1.2"
$ws.Cells.Item(33, 2).Value = 46270
$ws.Cells.Item(33, 4).Value = 100
Set-TextLiteral $ws.Cells.Item(33, 5) "2019-04-23"

$ws.Range("A27:E27").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)
$ws.Cells.Item(34, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.8.1"
$ws.Cells.Item(34, 2).Value = 48540
$ws.Cells.Item(34, 4).Value = 100
Set-TextLiteral $ws.Cells.Item(34, 5) "2019-03-16"

$ws.Range("A27:E27").Copy()
$ws.Range("A35:E35").PasteSpecial(-4122)
$ws.Cells.Item(35, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.5.2"
$ws.Cells.Item(35, 2).Value = 184564
$ws.Cells.Item(35, 4).Value = 99.90000000000001
Set-TextLiteral $ws.Cells.Item(35, 5) "2018-11-25"

$ws.Range("A27:E27").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)
$ws.Cells.Item(36, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.50.0.4"
$ws.Cells.Item(36, 2).Value = 14221
$ws.Cells.Item(36, 4).Value = 100
Set-TextLiteral $ws.Cells.Item(36, 5) "2018-05-08"

$ws.Range("A27:E27").Copy()
$ws.Range("A37:E37").PasteSpecial(-4122)
$ws.Cells.Item(37, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.30.1.2"
$ws.Cells.Item(37, 2).Value = 23765
$ws.Cells.Item(37, 4).Value = 100
Set-TextLiteral $ws.Cells.Item(37, 5) "2018-01-09"

$ws.Application.CutCopyMode = $false